$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Hospital names"
$ws.Range("A2").Value = "Radix Cosmo Dental"
$ws.Range("A3").Value = "Apollo Spectra Hospital"
$ws.Range("A4").Value = "Indraprastha Apollo Hospitals"
$ws.Range("A5").Value = "Nulife Hospital & Maternity Center"
$ws.Range("A6").Value = "Apollo Spectra Hospitals"
$ws.Range("A7").Value = "VIMHANS Nayati Superspecialty Hospital"
$ws.Range("A8").Value = "Apollo Cradle"
$ws.Range("A9").Value = "Handa Aesthetic and Plastics"
$ws.Range("A10").Value = "Fortis Escorts and Heart Institute"
$ws.Range("A11").Value = "Apollo Cradle"

$ws.Columns.Item(1).ColumnWidth = 36.6666666667
